# Apply the "add 2022-Q4 data" edit:
#  1. On the "总计" (summary) sheet, insert a new row for 2022-Q4 at the
#     top of the data (row 2), pushing every other quarter's row down by one.
#  2. Insert a brand-new worksheet named "2022-Q4" right after "总计",
#     holding the per-fund holdings detail for that quarter, and shift the
#     existing quarter sheets (2022-Q3 ... 2020-Q4) one slot to the right
#     (their own content is unchanged).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: "总计" summary sheet — shift rows 2..8 down to 3..9, then
# write the new 2022-Q4 row into row 2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Grab a style template for column A (bordered/bold "s=2" style) before
# we touch anything, so row 9 (brand new) can carry the same look as the
# other data rows in column A.
$styleSrc = $summary.Range("A8")
$styleSrc.Copy()
$summary.Range("A9").PasteSpecial(-4122)

# Shift existing data rows 8->9, 7->8, ... 2->3 (bottom-up so we never
# overwrite a row before we've read its old values). Column A carries a
# 0-based row index, so each value increments by one as it moves down.
$summary.Range("A9").Value = 7
$summary.Range("B9").Value = "2020-Q4"
$summary.Range("C9").Value = 3
$summary.Range("D9").Value = 0.13

$summary.Range("A8").Value = 6
$summary.Range("B8").Value = "2021-Q2"
$summary.Range("C8").Value = 10
$summary.Range("D8").Value = 1.8

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q3"
$summary.Range("C7").Value = 3
$summary.Range("D7").Value = 1.32

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q4"
$summary.Range("C6").Value = 5
$summary.Range("D6").Value = 1.13

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 5
$summary.Range("D5").Value = 1.06

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 6
$summary.Range("D4").Value = 1.13

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 7
$summary.Range("D3").Value = 0.85

# New row 2 -- the 2022-Q4 summary line.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 0.79

Write-Host "summary sheet updated"

# ---------------------------------------------------------------------
# Part 2: brand-new "2022-Q4" worksheet, inserted right after "总计"
# (i.e. before the current 2nd sheet, "2022-Q3"). The existing quarter
# sheets are left untouched content-wise; they simply slide one tab to
# the right because the new sheet is spliced in before them.
# ---------------------------------------------------------------------
$oldQ3 = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($oldQ3)
$q4.Name = "2022-Q4"

# Borrow the header/column-A formatting ("s=2": bold, thin border,
# center/top alignment) from the existing "2022-Q3" sheet so the new
# sheet matches the look of its siblings. (Copy/PasteSpecial cell-by-cell
# -- multi-cell-range-to-multi-cell-range paste does not carry styling.)
foreach ($col in @("B", "C", "D", "E", "F", "G", "H")) {
    $oldQ3.Range("$col`1").Copy()
    $q4.Range("$col`1").PasteSpecial(-4122)
}
foreach ($row in 2..8) {
    $oldQ3.Range("A2").Copy()
    $q4.Range("A$row").PasteSpecial(-4122)
}

# Header row.
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# The text-ish numeric columns (fund code, scale, position %, etc.) are
# stored as text in the source data, so force Text format before writing
# them -- otherwise Excel normalises them into numeric cells.
$textCols = $q4.Range("B2:G8")
$textCols.NumberFormat = "@"

# Row 2 -- 515210 国泰中证钢铁ETF
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "515210"
$q4.Range("C2").Value = "国泰中证钢铁ETF"
$q4.Range("D2").Value = "13.53"
$q4.Range("E2").Value = "99.48"
$q4.Range("F2").Value = "2.88"
$q4.Range("G2").Value = "0.3897"
$q4.Range("H2").Value = 8

# Row 3 -- 502023 鹏华国证钢铁行业指数（LOF）A
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "502023"
$q4.Range("C3").Value = "鹏华国证钢铁行业指数（LOF）A"
$q4.Range("D3").Value = "8.97"
$q4.Range("E3").Value = "94.55"
$q4.Range("F3").Value = "2.60"
$q4.Range("G3").Value = "0.2332"
$q4.Range("H3").Value = 9

# Row 4 -- 012810 鹏华国证钢铁行业指数（LOF）C
$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "012810"
$q4.Range("C4").Value = "鹏华国证钢铁行业指数（LOF）C"
$q4.Range("D4").Value = "3.28"
$q4.Range("E4").Value = "94.55"
$q4.Range("F4").Value = "2.60"
$q4.Range("G4").Value = "0.0853"
$q4.Range("H4").Value = 9

# Row 5 -- 168203 中融国证钢铁行业指数A
$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "168203"
$q4.Range("C5").Value = "中融国证钢铁行业指数A"
$q4.Range("D5").Value = "3.18"
$q4.Range("E5").Value = "92.38"
$q4.Range("F5").Value = "2.53"
$q4.Range("G5").Value = "0.0805"
$q4.Range("H5").Value = 9

# Row 6 -- 013802 财通资管中证钢铁指数A
$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "013802"
$q4.Range("C6").Value = "财通资管中证钢铁指数A"
$q4.Range("D6").Value = "0.08"
$q4.Range("E6").Value = "91.91"
$q4.Range("F6").Value = "2.73"
$q4.Range("G6").Value = "0.0022"
$q4.Range("H6").Value = 9

# Row 7 -- 013803 财通资管中证钢铁指数C
$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "013803"
$q4.Range("C7").Value = "财通资管中证钢铁指数C"
$q4.Range("D7").Value = "0.02"
$q4.Range("E7").Value = "91.91"
$q4.Range("F7").Value = "2.73"
$q4.Range("G7").Value = "0.0005"
$q4.Range("H7").Value = 9

# Row 8 -- 016815 中融国证钢铁行业指数C
$q4.Range("A8").Value = 6
$q4.Range("B8").Value = "016815"
$q4.Range("C8").Value = "中融国证钢铁行业指数C"
$q4.Range("D8").Value = "0.02"
$q4.Range("E8").Value = "92.38"
$q4.Range("F8").Value = "2.53"
$q4.Range("G8").Value = "0.0005"
$q4.Range("H8").Value = 9

Write-Host "2022-Q4 sheet created"
